# [服務商優惠徵集] Paid.xlsx — "Add files via upload"
#
# The WotoHub / WotoKOL (臥兔) row (row 7 of "Form Responses 1") is updated:
#   - A7: contact name corrected "Lucia Chang" -> "Lucia Zhang"
#   - E7: the benefit description is replaced with the final, longer
#     benefit copy, wrapped + vertically centered, with the row resized
#     to fit the new text.
# The previously selected cell (E2) is moved to E8 to reflect where the
# editor ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

# --- Row 7 (WotoHub / WotoKOL) -------------------------------------------

# Contact name correction.
$ws.Range("A7").Value = "Lucia Zhang"

# Updated benefit copy.
$ws.Range("E7").Value = "臥兔6000萬網紅紅人庫 專業版提供每月5000封網紅郵件聯繫額度, 加入亞馬遜賣家成長服務 加送每月額外5000封網紅郵件聯繫額度"

# Match the formatting used elsewhere for the longer benefit cells:
# wrap text + vertically centered.
$ws.Range("E7").WrapText = $true
$ws.Range("E7").VerticalAlignment = -4108  # xlCenter

# Row grows to fit the longer text.
$ws.Rows(7).RowHeight = 90.5

# --- Selection -------------------------------------------------------------
# Leave the selection where the editor ended up.
$ws.Range("E8").Select()
